$wb = $excel.ActiveWorkbook
$nl = [char]10

$ws1 = $wb.Worksheets.Item(1)   # Pediatric VFC Vaccine
$ws2 = $wb.Worksheets.Item(2)   # Adult Vaccine
$ws3 = $wb.Worksheets.Item(3)   # Pediatric Influenza Vaccine
$ws4 = $wb.Worksheets.Item(4)   # Adult Influenza Vaccine

# ---- Sheet 1: Pediatric VFC Vaccine ----
$r = $ws1.Cells.Replace("DTaP [1]", "DTaP ")
$r = $ws1.Cells.Replace("DTaP-IPV [2]", "DTaP-IPV ")
$r = $ws1.Cells.Replace("DTaP-Hep B-IPV [4]", "DTaP-Hep B-IPV ")
$r = $ws1.Cells.Replace("DTaP-IP-HI [4]", "DTaP-IP-HI ")
$r = $ws1.Cells.Replace("e-IPV [5]", "e-IPV ")
$r = $ws1.Cells.Replace("Hepatitis A Pediatric [5]", "Hepatitis A Pediatric ")
$r = $ws1.Cells.Replace("Hepatitis A-Hepatitis B 18 only [3]", "Hepatitis A-Hepatitis B 18 only ")

$old = "Hepatitis B [5]" + $nl + "Pediatric/Adolescent"
$r = $ws1.Cells.Replace($old, "Hepatitis B  Pediatric/Adolescent")

$old = "Recombivax" + $nl + "HB"
$r = $ws1.Cells.Replace($old, "Recombivax HB")

$r = $ws1.Cells.Replace("Hib [5]", "Hib ")
$r = $ws1.Cells.Replace("HPV - Human Papillomavirus 9-valent [5]", "HPV - Human Papillomavirus 9-valent ")
$r = $ws1.Cells.Replace("MENB - Meningococcal Group B [5]", "MENB - Meningococcal Group B ")
$r = $ws1.Cells.Replace("Meningococcal Conjugate (Groups A, C, Y and W-135) [5]", "Meningococcal Conjugate (Groups A, C, Y and W-135) ")
$r = $ws1.Cells.Replace("Measles, Mumps and Rubella (MMR) [1]", "Measles, Mumps and Rubella (MMR) ")
$r = $ws1.Cells.Replace("MMR/Varicella [2]", "MMR/Varicella ")

$old = "Pneumococcal" + $nl + "13-valent [5] (Pediatric)"
$r = $ws1.Cells.Replace($old, "Pneumococcal 13-valent  (Pediatric)")

$r = $ws1.Cells.Replace("Rotavirus, Live, Oral, Pentavalent [5]", "Rotavirus, Live, Oral, Pentavalent ")
$r = $ws1.Cells.Replace("Rotavirus, Live, Oral, Oral [5]", "Rotavirus, Live, Oral, Oral ")
$r = $ws1.Cells.Replace("Tetanus and Diphtheria Toxoids [3]", "Tetanus and Diphtheria Toxoids ")
$r = $ws1.Cells.Replace("Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]", "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis ")
$r = $ws1.Cells.Replace("Varicella [5]", "Varicella ")

# ---- Sheet 2: Adult Vaccine ----
$r = $ws2.Cells.Replace("Hepatitis A Adult [5]", "Hepatitis A Adult ")
$r = $ws2.Cells.Replace("Hepatitis A-Hepatitis B Adult [3]", "Hepatitis A-Hepatitis B Adult ")
$r = $ws2.Cells.Replace("Hepatitis B Adult [5]", "Hepatitis B Adult ")
$r = $ws2.Cells.Replace("HPV-Human Papillomavirus 9 Valent [5]", "HPV-Human Papillomavirus 9 Valent ")
$r = $ws2.Cells.Replace("Measles, Mumps,  Rubella [1]", "Measles, Mumps,  Rubella ")
$r = $ws2.Cells.Replace("Meningococcal Conjugate (Groups A, C, Y and W-135) [5]", "Meningococcal Conjugate (Groups A, C, Y and W-135) ")
$r = $ws2.Cells.Replace("MENB - Meningococcal Group B [5]", "MENB - Meningococcal Group B ")

$old = "Pneumococcal" + $nl + "13-valent [5]"
$r = $ws2.Cells.Replace($old, "Pneumococcal 13-valent ")

$r = $ws2.Cells.Replace("Tetanus and Diphtheria Toxoids [3]", "Tetanus and Diphtheria Toxoids ")
$r = $ws2.Cells.Replace("Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis [1]", "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis ")
$r = $ws2.Cells.Replace("Varicella [5]", "Varicella ")

# ---- Sheet 3: Pediatric Influenza Vaccine ----
$old = "Influenza [5]" + $nl + "(Age 6 months and older)"
$r = $ws3.Cells.Replace($old, "Influenza  (Age 6 months and older)")

$old = "Fluzone" + $nl + "Quadrivalent" + $nl + "Pediatric dose"
$r = $ws3.Cells.Replace($old, "Fluzone Quadrivalent Pediatric dose")

$old = "Fluzone" + $nl + "Quadrivalent"
$r = $ws3.Cells.Replace($old, "Fluzone Quadrivalent")

$old = "Influenza [5]" + $nl + "(Age 6-35 months)"
$r = $ws3.Cells.Replace($old, "Influenza  (Age 6-35 months)")

$old = "Fluarix" + $nl + "Quadrivalent"
$r = $ws3.Cells.Replace($old, "Fluarix Quadrivalent")

$old = "FluLaval" + $nl + "Quadrivalent"
$r = $ws3.Cells.Replace($old, "FluLaval Quadrivalent")

$old = "Influenza [5]" + $nl + "(Age 4 years and older)"
$r = $ws3.Cells.Replace($old, "Influenza  (Age 4 years and older)")

$old = "Influenza [5]" + $nl + "(Age 6 -35 months)"
$r = $ws3.Cells.Replace($old, "Influenza  (Age 6 -35 months)")

$old = "Influenza [5]" + $nl + "(Age 36 months and older)"
$r = $ws3.Cells.Replace($old, "Influenza  (Age 36 months and older)")

$old = "Influenza [5]" + $nl + "Live, Intranasal (Age 2-49 years)"
$r = $ws3.Cells.Replace($old, "Influenza  Live, Intranasal (Age 2-49 years)")

$old = "FluMist" + $nl + "Quadrivalent"
$r = $ws3.Cells.Replace($old, "FluMist Quadrivalent")

# ---- Sheet 4: Adult Influenza Vaccine ----
$old = "Influenza [5]" + $nl + "(Age 6 months and older)"
$r = $ws4.Cells.Replace($old, "Influenza  (Age 6 months and older)")

$old = "Fluzone" + $nl + "Quadrivalent"
$r = $ws4.Cells.Replace($old, "Fluzone Quadrivalent")

$old = "Fluarix" + $nl + "Quadrivalent"
$r = $ws4.Cells.Replace($old, "Fluarix Quadrivalent")

$old = "FluLaval" + $nl + "Quadrivalent"
$r = $ws4.Cells.Replace($old, "FluLaval Quadrivalent")

$old = "Influenza [5]" + $nl + "(Age 4 years and older)"
$r = $ws4.Cells.Replace($old, "Influenza  (Age 4 years and older)")

$old = "Influenza [5]" + $nl + "(Age 36 months and older)"
$r = $ws4.Cells.Replace($old, "Influenza  (Age 36 months and older)")

$old = "Afluria" + $nl + "Quadrivalent"
$r = $ws4.Cells.Replace($old, "Afluria Quadrivalent")

Write-Host "All replacements done."
